$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 41: additional receiver VR2W-137168 / Denton Slough
$ws.Range("A41").Value = "VR2W-137168"
$ws.Range("D41").Value = "Denton Slough"
$ws.Range("B41").Value = 48.203850000000003
$ws.Range("C41").Value = -116.25677

# Row 6: receiver id changed (VR2W-134451 -> VR2W-137164)
$ws.Range("A6").Value = "VR2W-137164"

# Row 21: site name renamed (Long Bridge -> Long Bridge S.)
$ws.Range("D21").Value = "Long Bridge S."

# Row 18: site name renamed (Long Bridge -> Long Bridge N.)
$ws.Range("D18").Value = "Long Bridge N."

# Row 38: site name renamed (Kootenai Point -> Kootenai Bay)
$ws.Range("D38").Value = "Kootenai Bay"

# Row 27: site name spelling fix (Mouth of Morton Sough -> Mouth of Morton Slough)
$ws.Range("D27").Value = "Mouth of Morton Slough"

# Row 39: updated coordinates for VR2W-137167 / Oden Bay
$ws.Range("B39").Value = 48.305210000000002
$ws.Range("C39").Value = -116.45361

# Page setup: portrait orientation with explicit print resolution (300 dpi)
$ws.PageSetup.Orientation = 1
